$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Typist" and "Typist QC" columns (E:F) are no longer needed -
# remove them entirely so every column to their right shifts left by two.
$ws.Range("E1:F1").EntireColumn.Delete()

# The "Process" value for the first data row was missing; the Typist/Typist QC
# columns are being replaced conceptually by a single "Production & QC" status,
# recorded against the (now shifted) Process column for that order.
$ws.Range("G2").Value = "Production & QC"

# Re-fit the columns whose contents changed so the sheet reads cleanly
# (match the width Excel's own best-fit would have produced for the new text).
$ws.Range("B1").EntireColumn.ColumnWidth = 9
$ws.Range("E1").EntireColumn.ColumnWidth = 19.833333333333332
$ws.Range("G1").EntireColumn.ColumnWidth = 13.5

# The old conditional formatting rule was scoped to the (now removed) Typist
# column, so drop it rather than let it silently point at the wrong cell.
$ws.Cells.FormatConditions.Delete()

# Reset the view's selection away from the old Typist/Typist QC area.
$ws.Range("A4:XFD5").Select()
